# Update cryptocurrency price (D) and 1h volume-change (E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain numeric-looking text (e.g. "591.24", "33.00",
# "3.300.50") that must stay TEXT, not be reinterpreted as numbers (which
# would round/reformat them). Force text number-format before assigning,
# then restore the "Normal" style so no visible formatting changes.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.253.16'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.607.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.82%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '190.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.36%  '
$ws.Range("E7").Value = '  +1.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.601.05'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.82%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  +3.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.663'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.68'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.69%  '
$ws.Range("E13").Value = '  +3.89%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.88'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.186.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.72'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.607.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.225.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.32%  '
$ws.Range("E19").Value = '  +4.37%  '
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("E21").Value = '  +4.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '492.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.10'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.15%  '
$ws.Range("E24").Value = '  +6.82%  '
$ws.Range("E25").Value = '  +4.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.74%  '
$ws.Range("E28").Value = '  +1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '637.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.48%  '
$ws.Range("E34").Value = '  +6.65%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.72'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("E36").Value = '  +6.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0818'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.405'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.27%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.55'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.300.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0452'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.43%  '
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("E47").Value = '  +2.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.32%  '
$ws.Range("E49").Value = '  -2.26%  '
$ws.Range("E50").Value = '  +4.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.03%  '
